$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last four forecast-horizon rows (Q4..Q7) — the simulated rt_data
# for those components was bugged, so the evaluation table now only keeps
# rows for Q0..Q3 (rows 2-5). Deleting the rows also drops the now-unused
# "Q4".."Q7" shared strings automatically.
$ws.Range("A6:G9").EntireRow.Delete() | Out-Null

# Refresh the forecast-error figures for the remaining rows with the
# corrected values.
$ws.Range("B2").Value = -0.2902542872142831
$ws.Range("C2").Value = 0.4474900817361138
$ws.Range("D2").Value = 0.3638477001211933
$ws.Range("E2").Value = 0.6031978946591187
$ws.Range("F2").Value = 0.5487329520217746

$ws.Range("B3").Value = -0.08670396990747065
$ws.Range("C3").Value = 0.4011283099121256
$ws.Range("D3").Value = 0.2437355075710951
$ws.Range("E3").Value = 0.4936957641818442
$ws.Range("F3").Value = 0.5123127616487798
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = -0.2842180920556761
$ws.Range("C4").Value = 0.4543513181817999
$ws.Range("D4").Value = 0.2495528846253024
$ws.Range("E4").Value = 0.4995526845341764
$ws.Range("F4").Value = 0.4500306133234053
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = -0.2239409044200031
$ws.Range("C5").Value = 0.5803388081800324
$ws.Range("D5").Value = 0.3869426609522694
$ws.Range("E5").Value = 0.6220471533189983
$ws.Range("F5").Value = 0.8207230132996398
$ws.Range("G5").Value = 2
